$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 298
$ws.Range("C2:C$lastRow").Value = 45205
